$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Section 1: Tests for shorts and connections ---
# Populate the "Observed" resistance measurements in column K for rows 18-33
$ws.Range("K18").Value = 0.8
$ws.Range("K19").Value = 0.7
$ws.Range("K20").Value = 0.7
$ws.Range("K21").Value = 0.7
$ws.Range("K22").Value = 0.7
$ws.Range("K23").Value = 0.7
$ws.Range("K24").Value = 0.7
$ws.Range("K25").Value = 0.7
$ws.Range("K26").Value = 0.7
$ws.Range("K27").Value = 0.7
$ws.Range("K28").Value = 0.8
$ws.Range("K29").Value = 0.7
$ws.Range("K30").Value = 0.7
$ws.Range("K31").Value = 0.7
$ws.Range("K32").Value = 0.8
$ws.Range("K33").Value = 0.7

# Observed resistance measurements in column K for rows 40-47
$ws.Range("K40").Value = 0.3
$ws.Range("K41").Value = 0.3
$ws.Range("K42").Value = 0
$ws.Range("K43").Value = 0.3
$ws.Range("K44").Value = 0.3
$ws.Range("K45").Value = 0
$ws.Range("K46").Value = 0.3
$ws.Range("K47").Value = 0.3

# --- Section 2: LED test ---
# Visual inspection answers
$ws.Range("M57").Value = "y"
$ws.Range("M58").Value = "y"
$ws.Range("M59").Value = "y"
# Result
$ws.Range("B60").Value = "pass"

# --- Section 3: HV test ---
# Voltage / current readings
$ws.Range("K66").Value = 98.8
$ws.Range("K67").Value = 9.877

# Voltage drop readings (no adapter / with adapter)
$ws.Range("C70").Value = -0.0174
$ws.Range("C71").Value = -0.017

# --- Sign off ---
$ws.Range("B80").Value = "Amanda"
$ws.Range("F80").Value = 41918

# --- Update selection to reflect where the tester left off ---
$ws.Range("F81").Select()
